# Applies the cell-value updates for the "Updated cryptos list" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every new value is written with a leading literal apostrophe so Excel
# stores it as text (matching the original inlineStr cells) instead of
# auto-converting numeric-looking strings (e.g. "19.61") into numbers.
$ws.Range("D2").Value = "'26.652.10"
$ws.Range("D3").Value = "'1.596.91"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("E4").Value = "'  +0.42%  "
$ws.Range("E5").Value = "'  -0.20%  "
$ws.Range("E6").Value = "'  -0.33%  "
$ws.Range("E7").Value = "'  +0.44%  "
$ws.Range("E8").Value = "'  -0.44%  "
$ws.Range("D10").Value = "'19.61"
$ws.Range("E10").Value = "'  +0.23%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "'  +0.46%  "
$ws.Range("D12").Value = "'1.821.46"
$ws.Range("E12").Value = "'  -0.17%  "
$ws.Range("D13").Value = "'1.599.20"
$ws.Range("E13").Value = "'  +0.24%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "'  +0.04%  "
$ws.Range("E16").Value = "'  -0.87%  "
$ws.Range("D17").Value = "'26.633.44"
$ws.Range("E17").Value = "'  -0.16%  "
$ws.Range("E18").Value = "'  -2.86%  "
$ws.Range("E19").Value = "'  +0.35%  "
$ws.Range("D20").Value = "'208.49"
$ws.Range("E20").Value = "'  -0.82%  "
$ws.Range("D21").Value = "'7.12"
$ws.Range("E21").Value = "'  -1.19%  "
$ws.Range("E22").Value = "'  -0.04%  "
$ws.Range("E23").Value = "'  -2.70%  "
$ws.Range("E24").Value = "'  +0.09%  "
$ws.Range("D25").Value = "'143.89"
$ws.Range("E25").Value = "'  +0.54%  "
$ws.Range("E26").Value = "'  +0.49%  "
$ws.Range("D27").Value = "'7.15"
$ws.Range("E27").Value = "'  +0.26%  "
$ws.Range("E28").Value = "'  -0.67%  "
$ws.Range("D29").Value = "'15.28"
$ws.Range("E29").Value = "'  -0.28%  "
$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "'  -2.46%  "
$ws.Range("E31").Value = "'  -0.30%  "
$ws.Range("E32").Value = "'  -0.45%  "
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("E34").Value = "'  +19.12%  "
$ws.Range("D35").Value = "'1.279.36"
$ws.Range("E35").Value = "'  -0.89%  "
$ws.Range("D36").Value = "'2.50"
$ws.Range("E36").Value = "'  +1.13%  "
$ws.Range("D37").Value = "'0.600"
$ws.Range("E37").Value = "'  -3.46%  "
$ws.Range("E38").Value = "'  -0.90%  "
$ws.Range("E39").Value = "'  -2.03%  "
$ws.Range("D40").Value = "'0.823"
$ws.Range("E40").Value = "'  -0.48%  "
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "'  +1.80%  "
$ws.Range("E42").Value = "'  -0.22%  "
$ws.Range("D43").Value = "'0.776"
$ws.Range("E43").Value = "'  -1.08%  "
$ws.Range("D44").Value = "'62.56"
$ws.Range("E44").Value = "'  -0.96%  "
$ws.Range("D45").Value = "'1.733.00"
$ws.Range("E45").Value = "'  -0.09%  "
$ws.Range("D46").Value = "'90.27"
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "'  -0.38%  "
$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.0₆0104"
$ws.Range("E48").Value = "'  -2.31%  "
$ws.Range("B49").Value = "'Algorand"
$ws.Range("C49").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.103"
$ws.Range("E49").Value = "'  +2.19%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0512"
$ws.Range("E50").Value = "'  +0.58%  "
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.49"
$ws.Range("E51").Value = "'  +1.33%  "
